$wb = $excel.ActiveWorkbook

# --- "About" sheet (sheet 1) ---
$ws1 = $wb.Worksheets.Item(1)

# Drop the old hyperlink to the source article before removing the rows
# that held the old citation (page/author/title/url), per the diff.
$ws1.Hyperlinks.Delete()

# Remove the old "Source" detail rows (year, title, author, URL) - rows 4-7.
$ws1.Rows("4:7").Delete()

# Source now reads "None".
$ws1.Range("B3").Value = "None"

# New note explaining the US-specific override.
$ws1.Range("A9").Value = "In the US, we set this to 0 so that increasing EV chargers does not induce additional deployment."

# --- "EoCSoEVMS" sheet (sheet 2) ---
$ws2 = $wb.Worksheets.Item(2)

# Value set to 0 in the US so extra chargers don't induce extra EV deployment.
$ws2.Range("B2").Value = 0

# Remove the now-unused "Hyperlink" cell style.
$wb.Styles("Hyperlink").Delete()

# Restore the selections/active sheet to match the saved view state.
$ws2.Range("B3").Select()
$ws1.Range("A4:XFD7").Select()
